$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the numeric values in B2:B5 with their word equivalents
$ws.Range("B2").Value = "one"
$ws.Range("B3").Value = "two"
$ws.Range("B4").Value = "three"
$ws.Range("B5").Value = "four"

# Update the selected cell to match the saved selection in the workbook
$ws.Range("B5").Select()
